# Workbook2.xlsx — "More asserted_distribution test cases."
#
# The sheet gains three new blank spacer rows (before the existing rows 4,
# 5 and 13 of the original layout), shifting the existing data down so the
# sheet grows from 14 to 17 data rows. Using Rows.Insert() (rather than
# writing cell-by-cell) lets Excel copy the row-above formatting forward
# automatically, which matches the target styles (including the D-column
# style="3" that bleeds into the new row 15 from the row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert in ascending order of *final* row number so each subsequent
# insert targets the already-shifted sheet.
$ws.Rows.Item(4).Insert()   # new blank row 4  (old row 4 -> 5)
$ws.Rows.Item(6).Insert()   # new blank row 6  (old row 5 -> 7)
$ws.Rows.Item(15).Insert()  # new blank row 15 (old row 13 -> 16)

# Update the window view to match: scrolled down so row 13 is at the top,
# with the new blank row 15 selected (whole row).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(15).EntireRow.Select()
